# M2SL Price importer TV data pull update
# - Revise a run of recent M2SL observations (rows 253-281) with corrected values
# - Append a new observation row (282) for 2023-05-01
# - Refresh the SeriesInfo metadata sheet (realtime dates, observation_end,
#   last_updated, popularity) and append a new "Source" = "fred" row

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: Closing_Price (M2SL observations)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Closing_Price")

$revisions = @{
    253 = 19114.3
    254 = 19356.7
    255 = 19600.3
    256 = 19841.2
    257 = 20116.9
    258 = 20431.1
    259 = 20506.6
    260 = 20662.9
    261 = 20847.8
    262 = 20964.3
    263 = 21116.2
    264 = 21316.1
    265 = 21549.1
    266 = 21561.8
    267 = 21570.3
    268 = 21697.5
    269 = 21677.1
    270 = 21665.2
    271 = 21665.7
    272 = 21702.6
    273 = 21658.9
    275 = 21432.3
    276 = 21399
    278 = 21212.7
    279 = 21077.4
    280 = 20841.1
    281 = 20674.6
}

foreach ($r in $revisions.Keys) {
    $ws1.Cells.Item($r, 2).Value = $revisions[$r]
}

# Append the new observation row (copy formatting from the prior row, then
# overwrite with the new date/value so style s="3" carries over to A282)
$ws1.Range("A281").Copy($ws1.Range("A282"))
$ws1.Cells.Item(282, 1).Value = 45047
$ws1.Cells.Item(282, 2).Value = 20805.5

# ---------------------------------------------------------------------------
# Sheet 2: SeriesInfo (metadata)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("SeriesInfo")

# B3/B4/B7 hold plain ISO date TEXT (not real Excel dates) in the original
# file. Assigning an ISO-looking string straight to .Value lets Excel's
# auto-detection silently convert it into a date serial, so force the
# cells to Text format first and restore the (unstyled) Normal style
# afterwards to keep them looking exactly like the untouched neighbours.
$datesRng = $ws2.Range("B3:B7")
$datesRng.NumberFormat = "@"
$ws2.Cells.Item(3, 2).Value = "2023-07-13"
$ws2.Cells.Item(4, 2).Value = "2023-07-13"
$ws2.Cells.Item(7, 2).Value = "2023-05-01"
$datesRng.Style = "Normal"

$ws2.Cells.Item(14, 2).Value = "2023-06-27 12:03:02-05"
$ws2.Cells.Item(15, 2).Value = 93

# Append the new "Source" row (copy formatting from the prior row so A17
# keeps the bold "label" style s="1")
$ws2.Range("A16").Copy($ws2.Range("A17"))
$ws2.Cells.Item(17, 1).Value = "Source"
$ws2.Cells.Item(17, 2).Value = "fred"
